# This script reproduces the commit:
#   [expression] JSON >> keys(jsonpath): extract immediate keys of resolved
#   JSON fragment based on jsonpath.
#   [json] storeKeys(json,jsonpath,var): extract immediate keys of resolved
#   JSON fragment based on jsonpath.
#
# The macro catalogue lives on the hidden "#system" worksheet: row 1 holds
# the category name for each column, and each column below it lists the
# macro signatures that belong to that category (alphabetically). Named
# ranges ("target", "json", "web", ...) point at those columns.
#
# Two structural edits happened on that sheet:
#   1. The whole "text" category (column Y, a single-entry category) was
#      retired, so its column was deleted outright (Z:AE shift left to
#      Y:AD) and its entry was removed from the "target" master list
#      (column A), which itself shifts up by one starting where "text"
#      used to be listed.
#   2. A new JSON macro "storeKeys(json,jsonpath,var)" was added, sorted
#      alphabetically into column M, between "storeCount(...)" and
#      "storeValue(...)" - i.e. inserted at M16, pushing the two rows
#      below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------
# 1) Remove the retired "text" category column (Y) completely. Everything
#    to its right (web, webalert, webcookie, ws, xml, ws.async) shifts
#    one column to the left.
# ---------------------------------------------------------------------
$ws.Columns.Item("Y").Delete()
Write-Output "Deleted retired 'text' column (Y)."

# ---------------------------------------------------------------------
# 2) Remove "text" from the master category list in column A (it was
#    listed right after "step", before "web") and pull the remaining
#    entries (web, webalert, webcookie, ws, xml) up by one row.
# ---------------------------------------------------------------------
for ($r = 25; $r -le 30; $r++) {
    $ws.Cells.Item($r, 1).Value2 = $ws.Cells.Item($r + 1, 1).Value2
}
$ws.Cells.Item(31, 1).ClearContents()

# ---------------------------------------------------------------------
# 3) Insert the new "storeKeys(json,jsonpath,var)" macro into column M
#    (json category), keeping the alphabetical ordering: it belongs
#    right before "storeValue(json,jsonpath,var)" at row 16, so the
#    existing M16/M17 values shift down to M17/M18.
# ---------------------------------------------------------------------
$ws.Cells.Item(18, 13).Value2 = $ws.Cells.Item(17, 13).Value2
$ws.Cells.Item(17, 13).Value2 = $ws.Cells.Item(16, 13).Value2
$ws.Cells.Item(16, 13).Value2 = "storeKeys(json,jsonpath,var)"

# ---------------------------------------------------------------------
# 4) Fix up the named ranges that describe each category's extent. The
#    column delete above does not auto-repair these, and the "json"
#    range grew by one row because of the newly inserted macro.
# ---------------------------------------------------------------------
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$30"
$wb.Names.Item("json").RefersTo = "='#system'!`$M`$2:`$M`$18"
$wb.Names.Item("web").RefersTo = "='#system'!`$Y`$2:`$Y`$129"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AA`$2:`$AA`$8"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AB`$2:`$AB`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AC`$2:`$AC`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AD`$2:`$AD`$27"

Write-Output "Done: added storeKeys(json,jsonpath,var) and dropped the 'text' column."
